# Edits the "R Coding presentation" deck:
#   Slide 2 (Content Placeholder): append " age" to the "...based on" run
#   Slide 3 (Content Placeholder): four sentence rewrites about the data-wrangling steps
#   Slide 9 (Content Placeholder): append three new bullet paragraphs (one with text,
#                                   two intentionally blank) after the last bullet

$p = $ppt.ActivePresentation

function Replace-SubText {
    param(
        $TextRange,
        [string]$OldText,
        [string]$NewText
    )
    $full = $TextRange.Text
    $idx = $full.IndexOf($OldText)
    if ($idx -lt 0) {
        throw "Text not found: $OldText"
    }
    $start = $idx + 1
    $len = $OldText.Length
    $sub = $TextRange.Characters($start, $len)
    $sub.Text = $NewText
}

# --- Slide 2 : "What are we doing here?" -------------------------------
$slide2 = $p.Slides.Item(2)
$tr2 = $slide2.Shapes.Item(2).TextFrame.TextRange
Replace-SubText $tr2 ":  See if there are any interesting differences in the amygdala based on" ":  See if there are any interesting differences in the amygdala based on age"

# --- Slide 3 : "Overview of code" ---------------------------------------
$slide3 = $p.Slides.Item(3)
$tr3 = $slide3.Shapes.Item(2).TextFrame.TextRange
Replace-SubText $tr3 "Make heat map based on those age ranges" "Make heat map based on those age range averages"
Replace-SubText $tr3 "Sort by age range (in decades)" "Take the average of each gene in each age range"
Replace-SubText $tr3 "Combine amygdala and age data" "Combine amygdala and annotations data"
Replace-SubText $tr3 "Includes function to automatically create comment indicators rather than by hand" "Includes commenting out shape data of the given matrices"

# --- Slide 9 : "The future" ----------------------------------------------
$slide9 = $p.Slides.Item(9)
$tr9 = $slide9.Shapes.Item(2).TextFrame.TextRange
$nl = [char]13
$tr9.InsertAfter($nl + "Isolate genes with high variation across age ranges") | Out-Null
$tr9.InsertAfter($nl) | Out-Null
$tr9.InsertAfter($nl) | Out-Null
